# Append a new "Messi vs Ronaldo Stats:" section (with its source hyperlink)
# right after the existing "FIFA quote:" hyperlink paragraph, keeping the
# document's final trailing empty paragraph intact.

$d = $word.ActiveDocument

# Locate the paragraph that holds the FIFA-quote hyperlink by searching for
# its well-known URL text instead of a hard-coded paragraph index, so the
# script is resilient to the exact paragraph numbering of the document.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*goal.com*fifa-player-ratings-explained*") {
        $target = $para
    }
}

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Build the three new paragraphs as a single OOXML fragment:
#   1. a blank separator paragraph
#   2. "Messi vs Ronaldo Stats:" split across three runs (matches the
#      original commit's run layout)
#   3. a blank paragraph that will host the new hyperlink
$block = "<w:p $ns/>" +
    "<w:p $ns><w:r><w:t xml:space='preserve'>Messi vs </w:t></w:r><w:r><w:t>Ronaldo</w:t></w:r><w:r><w:t xml:space='preserve'> Stats:</w:t></w:r></w:p>" +
    "<w:p $ns/>"

$insertAt = $target.Range.End
$rng = $d.Range($insertAt, $insertAt)
$rng.InsertXML($block)

# The third inserted paragraph (now immediately after the "Messi vs
# Ronaldo Stats:" paragraph) is still empty; turn it into a real hyperlink
# using the high-level API so Word wires up the relationship / styling
# exactly like it does for the document's other hyperlinks.
$hyperlinkPara = $target.Next().Next().Next()
$hStart = $hyperlinkPara.Range.Start
$hRange = $d.Range($hStart, $hStart)
$d.Hyperlinks.Add($hRange, "https://www.messivsronaldo.app/calendar-year-stats/2020") | Out-Null
